$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 439, shifting existing rows 439:466 down to 440:467.
$ws.Rows(439).Insert()

# Populate the newly inserted row 439 with its data.
$ws.Cells.Item(439, 1).Value = 8
$ws.Cells.Item(439, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(439, 3).Value = "Coquimbo"
$ws.Cells.Item(439, 4).Value = 44610
$ws.Cells.Item(439, 5).Value = 4
$ws.Cells.Item(439, 6).Value = 100112043
$ws.Cells.Item(439, 7).Value = "Pepino ensalada"
$ws.Cells.Item(439, 8).Value = "Sin especificar"
$ws.Cells.Item(439, 9).Value = "Primera"
$ws.Cells.Item(439, 10).Value = 720
$ws.Cells.Item(439, 11).Value = 13000
$ws.Cells.Item(439, 12).Value = 14000
$ws.Cells.Item(439, 13).Value = 13500
$ws.Cells.Item(439, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(439, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(439, 16).Value = 225
$ws.Cells.Item(439, 17).Value = 60
$ws.Cells.Item(439, 18).Value = "Hortaliza"
